$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.109.65'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '2.264.79'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '151.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15,046.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '93.80'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.485'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '33.12'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0801'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '2.616.88'
$ws.Range("E15").Value = '  -0.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").Value = '2.266.80'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.784'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.48%  '
$ws.Range("D19").Value = '41.987.57'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.71%  '
$ws.Range("D21").Value = '0.0₃0914'
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.06'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '243.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.66%  '
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.09'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0744'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.98%  '
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("E40").Value = '  +0.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '19.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("D44").Value = '2.007.59'
$ws.Range("E44").Value = '  -3.15%  '
$ws.Range("E45").Value = '  +11.17%  '
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.95%  '
$ws.Range("E51").Value = '  +0.55%  '
